$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 25.35940266666667
$ws.Cells.Item(2, 8).Value = 76.078208
$ws.Cells.Item(2, 9).Value = 0.005186643687654987
$ws.Cells.Item(2, 10).Value = 0.005186643687654986
$ws.Cells.Item(2, 13).Value = 16.27546433333333
$ws.Cells.Item(2, 14).Value = 48.826393
$ws.Cells.Item(2, 15).Value = 0.06628560529319844
$ws.Cells.Item(2, 16).Value = 0.06628560529319844
$ws.Cells.Item(2, 17).Value = 412.7360536159715
$ws.Cells.Item(2, 18).Value = 3714.624482543744
$ws.Cells.Item(2, 19).Value = 0.0003437998162763576
$ws.Cells.Item(2, 20).Value = 0.0003437998162763576
$ws.Cells.Item(3, 7).Value = 25.35940266666667
$ws.Cells.Item(3, 8).Value = 76.078208
$ws.Cells.Item(3, 9).Value = 0.005186643687654987
$ws.Cells.Item(3, 10).Value = 0.005186643687654986
$ws.Cells.Item(3, 15).Value = 0.3480686258826592
$ws.Cells.Item(3, 16).Value = 0.3480686258826592
$ws.Cells.Item(3, 17).Value = 2167.295152528143
$ws.Cells.Item(3, 18).Value = 19505.65637275328
$ws.Cells.Item(3, 19).Value = 0.001805307941305039
$ws.Cells.Item(3, 20).Value = 0.001805307941305039
$ws.Cells.Item(4, 7).Value = 25.35940266666667
$ws.Cells.Item(4, 8).Value = 76.078208
$ws.Cells.Item(4, 9).Value = 0.005186643687654987
$ws.Cells.Item(4, 10).Value = 0.005186643687654986
$ws.Cells.Item(4, 13).Value = 42.61351133333333
$ws.Cells.Item(4, 14).Value = 127.840534
$ws.Cells.Item(4, 15).Value = 0.17355341356458
$ws.Cells.Item(4, 16).Value = 0.17355341356458
$ws.Cells.Item(4, 17).Value = 1080.653192942563
$ws.Cells.Item(4, 18).Value = 9725.878736483071
$ws.Cells.Item(4, 19).Value = 0.0009001597169357041
$ws.Cells.Item(4, 20).Value = 0.0009001597169357041
$ws.Cells.Item(5, 7).Value = 25.35940266666667
$ws.Cells.Item(5, 8).Value = 76.078208
$ws.Cells.Item(5, 9).Value = 0.005186643687654987
$ws.Cells.Item(5, 10).Value = 0.005186643687654986
$ws.Cells.Item(5, 13).Value = 101.183272
$ws.Cells.Item(5, 14).Value = 303.549816
$ws.Cells.Item(5, 15).Value = 0.4120923552595624
$ws.Cells.Item(5, 16).Value = 0.4120923552595624
$ws.Cells.Item(5, 17).Value = 2565.947337778859
$ws.Cells.Item(5, 18).Value = 23093.52604000972
$ws.Cells.Item(5, 19).Value = 0.002137376213137886
$ws.Cells.Item(5, 20).Value = 0.002137376213137885
$ws.Cells.Item(6, 9).Value = 0.9837462940761621
$ws.Cells.Item(6, 10).Value = 0.983746294076162
$ws.Cells.Item(6, 13).Value = 16.27546433333333
$ws.Cells.Item(6, 14).Value = 48.826393
$ws.Cells.Item(6, 15).Value = 0.06628560529319844
$ws.Cells.Item(6, 16).Value = 0.06628560529319844
$ws.Cells.Item(6, 17).Value = 78283.29602489187
$ws.Cells.Item(6, 18).Value = 704549.6642240268
$ws.Cells.Item(6, 19).Value = 0.0652082185577792
$ws.Cells.Item(6, 20).Value = 0.06520821855777918
$ws.Cells.Item(7, 9).Value = 0.9837462940761621
$ws.Cells.Item(7, 10).Value = 0.983746294076162
$ws.Cells.Item(7, 15).Value = 0.3480686258826592
$ws.Cells.Item(7, 16).Value = 0.3480686258826592
$ws.Cells.Item(7, 19).Value = 0.3424112207962481
$ws.Cells.Item(7, 20).Value = 0.342411220796248
$ws.Cells.Item(8, 9).Value = 0.9837462940761621
$ws.Cells.Item(8, 10).Value = 0.983746294076162
$ws.Cells.Item(8, 13).Value = 42.61351133333333
$ws.Cells.Item(8, 14).Value = 127.840534
$ws.Cells.Item(8, 15).Value = 0.17355341356458
$ws.Cells.Item(8, 16).Value = 0.17355341356458
$ws.Cells.Item(8, 17).Value = 204966.5714012963
$ws.Cells.Item(8, 18).Value = 1844699.142611667
$ws.Cells.Item(8, 19).Value = 0.1707325274184231
$ws.Cells.Item(8, 20).Value = 0.1707325274184231
$ws.Cells.Item(9, 9).Value = 0.9837462940761621
$ws.Cells.Item(9, 10).Value = 0.983746294076162
$ws.Cells.Item(9, 13).Value = 101.183272
$ws.Cells.Item(9, 14).Value = 303.549816
$ws.Cells.Item(9, 15).Value = 0.4120923552595624
$ws.Cells.Item(9, 16).Value = 0.4120923552595624
$ws.Cells.Item(9, 17).Value = 486681.0477732701
$ws.Cells.Item(9, 18).Value = 4380129.42995943
$ws.Cells.Item(9, 19).Value = 0.4053943273037117
$ws.Cells.Item(9, 20).Value = 0.4053943273037117
$ws.Cells.Item(10, 7).Value = 51.27300266666666
$ws.Cells.Item(10, 8).Value = 153.819008
$ws.Cells.Item(10, 9).Value = 0.01048663484403512
$ws.Cells.Item(10, 10).Value = 0.01048663484403512
$ws.Cells.Item(10, 13).Value = 16.27546433333333
$ws.Cells.Item(10, 14).Value = 48.826393
$ws.Cells.Item(10, 15).Value = 0.06628560529319844
$ws.Cells.Item(10, 16).Value = 0.06628560529319844
$ws.Cells.Item(10, 17).Value = 834.4919261642381
$ws.Cells.Item(10, 18).Value = 7510.427335478144
$ws.Cells.Item(10, 19).Value = 0.0006951129381256138
$ws.Cells.Item(10, 20).Value = 0.0006951129381256138
$ws.Cells.Item(11, 7).Value = 51.27300266666666
$ws.Cells.Item(11, 8).Value = 153.819008
$ws.Cells.Item(11, 9).Value = 0.01048663484403512
$ws.Cells.Item(11, 10).Value = 0.01048663484403512
$ws.Cells.Item(11, 15).Value = 0.3480686258826592
$ws.Cells.Item(11, 16).Value = 0.3480686258826592
$ws.Cells.Item(11, 17).Value = 4381.953770586809
$ws.Cells.Item(11, 18).Value = 39437.58393528128
$ws.Cells.Item(11, 19).Value = 0.00365006858029652
$ws.Cells.Item(11, 20).Value = 0.00365006858029652
$ws.Cells.Item(12, 7).Value = 51.27300266666666
$ws.Cells.Item(12, 8).Value = 153.819008
$ws.Cells.Item(12, 9).Value = 0.01048663484403512
$ws.Cells.Item(12, 10).Value = 0.01048663484403512
$ws.Cells.Item(12, 13).Value = 42.61351133333333
$ws.Cells.Item(12, 14).Value = 127.840534
$ws.Cells.Item(12, 15).Value = 0.17355341356458
$ws.Cells.Item(12, 16).Value = 0.17355341356458
$ws.Cells.Item(12, 17).Value = 2184.92268023003
$ws.Cells.Item(12, 18).Value = 19664.30412207027
$ws.Cells.Item(12, 19).Value = 0.001819991273987563
$ws.Cells.Item(12, 20).Value = 0.001819991273987563
$ws.Cells.Item(13, 7).Value = 51.27300266666666
$ws.Cells.Item(13, 8).Value = 153.819008
$ws.Cells.Item(13, 9).Value = 0.01048663484403512
$ws.Cells.Item(13, 10).Value = 0.01048663484403512
$ws.Cells.Item(13, 13).Value = 101.183272
$ws.Cells.Item(13, 14).Value = 303.549816
$ws.Cells.Item(13, 15).Value = 0.4120923552595624
$ws.Cells.Item(13, 16).Value = 0.4120923552595624
$ws.Cells.Item(13, 17).Value = 5187.970175078058
$ws.Cells.Item(13, 18).Value = 46691.73157570252
$ws.Cells.Item(13, 19).Value = 0.004321462051625429
$ws.Cells.Item(13, 20).Value = 0.004321462051625429
$ws.Cells.Item(14, 7).Value = 2.837922333333333
$ws.Cells.Item(14, 8).Value = 8.513767
$ws.Cells.Item(14, 9).Value = 0.0005804273921477663
$ws.Cells.Item(14, 10).Value = 0.0005804273921477662
$ws.Cells.Item(14, 13).Value = 16.27546433333333
$ws.Cells.Item(14, 14).Value = 48.826393
$ws.Cells.Item(14, 15).Value = 0.06628560529319844
$ws.Cells.Item(14, 16).Value = 0.06628560529319844
$ws.Cells.Item(14, 17).Value = 46.18850371693677
$ws.Cells.Item(14, 18).Value = 415.696533452431
$ws.Cells.Item(14, 19).Value = 0.00003847398101726734
$ws.Cells.Item(14, 20).Value = 0.00003847398101726734
$ws.Cells.Item(15, 7).Value = 2.837922333333333
$ws.Cells.Item(15, 8).Value = 8.513767
$ws.Cells.Item(15, 9).Value = 0.0005804273921477663
$ws.Cells.Item(15, 10).Value = 0.0005804273921477662
$ws.Cells.Item(15, 15).Value = 0.3480686258826592
$ws.Cells.Item(15, 16).Value = 0.3480686258826592
$ws.Cells.Item(15, 17).Value = 242.5378624698161
$ws.Cells.Item(15, 18).Value = 2182.840762228345
$ws.Cells.Item(15, 19).Value = 0.0002020285648095284
$ws.Cells.Item(15, 20).Value = 0.0002020285648095283
$ws.Cells.Item(16, 7).Value = 2.837922333333333
$ws.Cells.Item(16, 8).Value = 8.513767
$ws.Cells.Item(16, 9).Value = 0.0005804273921477663
$ws.Cells.Item(16, 10).Value = 0.0005804273921477662
$ws.Cells.Item(16, 13).Value = 42.61351133333333
$ws.Cells.Item(16, 14).Value = 127.840534
$ws.Cells.Item(16, 15).Value = 0.17355341356458
$ws.Cells.Item(16, 16).Value = 0.17355341356458
$ws.Cells.Item(16, 17).Value = 120.9338355146198
$ws.Cells.Item(16, 18).Value = 1088.404519631578
$ws.Cells.Item(16, 19).Value = 0.0001007351552336319
$ws.Cells.Item(16, 20).Value = 0.0001007351552336319
$ws.Cells.Item(17, 7).Value = 2.837922333333333
$ws.Cells.Item(17, 8).Value = 8.513767
$ws.Cells.Item(17, 9).Value = 0.0005804273921477663
$ws.Cells.Item(17, 10).Value = 0.0005804273921477662
$ws.Cells.Item(17, 13).Value = 101.183272
$ws.Cells.Item(17, 14).Value = 303.549816
$ws.Cells.Item(17, 15).Value = 0.4120923552595624
$ws.Cells.Item(17, 16).Value = 0.4120923552595624
$ws.Cells.Item(17, 17).Value = 287.1502673685413
$ws.Cells.Item(17, 18).Value = 2584.352406316872
$ws.Cells.Item(17, 19).Value = 0.0002391896910873386
$ws.Cells.Item(17, 20).Value = 0.0002391896910873386
